$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-15) holds a "Förändrad" (changed) date that needs to be
# updated from 2023-10-22 (serial 45221) to 2023-10-25 (serial 45224).
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45224
}
